# 23/2-2018 Player Implementation Unstable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new laboration-diary entry for 23/2-2018 on row 10:
#   Datum | Uppgift | Timmar | Minuter
$ws.Range("A10").Value = "23/2-2018"
$ws.Range("B10").Value = "Player Implementation"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0

# The "Total tid" row (row 13) contains formulas that sum C2:C12/D2:D12,
# so they will pick up the new row automatically on recalculation.
$wb.Application.Calculate()

# Reflect the author's last active selection when they saved the file.
$ws.Range("F12").Select() | Out-Null
